$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "57.031.21"
Set-TextValue $ws.Range("E2") "  +0.74%  "
Set-TextValue $ws.Range("D3") "2.403.19"
Set-TextValue $ws.Range("E3") "  +1.05%  "
Set-TextValue $ws.Range("E4") "  +0.02%  "
Set-TextValue $ws.Range("D5") "506.83"
Set-TextValue $ws.Range("E5") "  -0.90%  "
Set-TextValue $ws.Range("D6") "132.97"
Set-TextValue $ws.Range("E6") "  +1.60%  "
Set-TextValue $ws.Range("E7") "  +0.13%  "
Set-TextValue $ws.Range("D8") "0.553"
Set-TextValue $ws.Range("E8") "  -0.19%  "
Set-TextValue $ws.Range("D9") "2.413.67"
Set-TextValue $ws.Range("E9") "  +0.52%  "
Set-TextValue $ws.Range("D10") "0.0977"
Set-TextValue $ws.Range("E10") "  +1.39%  "
Set-TextValue $ws.Range("E11") "  -1.19%  "
Set-TextValue $ws.Range("E12") "  +0.97%  "
Set-TextValue $ws.Range("D13") "4.61"
Set-TextValue $ws.Range("E13") "  -2.88%  "
Set-TextValue $ws.Range("D14") "2.830.04"
Set-TextValue $ws.Range("E14") "  +1.06%  "
Set-TextValue $ws.Range("D15") "56.953.96"
Set-TextValue $ws.Range("E15") "  +0.82%  "
Set-TextValue $ws.Range("D16") "21.80"
Set-TextValue $ws.Range("E16") "  +0.53%  "
Set-TextValue $ws.Range("E17") "  +1.83%  "
Set-TextValue $ws.Range("D18") "2.445.63"
Set-TextValue $ws.Range("E18") "  +1.91%  "
Set-TextValue $ws.Range("E19") "  -0.38%  "
Set-TextValue $ws.Range("E20") "  -0.20%  "
Set-TextValue $ws.Range("D21") "310.92"
Set-TextValue $ws.Range("E21") "  -0.63%  "
Set-TextValue $ws.Range("D22") "6.27"
Set-TextValue $ws.Range("E22") "  +0.41%  "
Set-TextValue $ws.Range("E23") "  +0.02%  "
Set-TextValue $ws.Range("D24") "5.60"
Set-TextValue $ws.Range("E24") "  -4.34%  "
Set-TextValue $ws.Range("D25") "67.84"
Set-TextValue $ws.Range("E25") "  +3.79%  "
Set-TextValue $ws.Range("E26") "  +0.05%  "
Set-TextValue $ws.Range("D27") "0.376"
Set-TextValue $ws.Range("E27") "  -3.75%  "
Set-TextValue $ws.Range("E28") "  -0.62%  "
Set-TextValue $ws.Range("E29") "  +2.49%  "
Set-TextValue $ws.Range("D30") "176.08"
Set-TextValue $ws.Range("E30") "  +0.76%  "
Set-TextValue $ws.Range("E31") "  +1.03%  "
Set-TextValue $ws.Range("E32") "  -0.56%  "
Set-TextValue $ws.Range("D33") "1.12"
Set-TextValue $ws.Range("E33") "  +0.70%  "
Set-TextValue $ws.Range("D34") "5.90"
Set-TextValue $ws.Range("E34") "  -4.46%  "
Set-TextValue $ws.Range("E35") "  +0.23%  "
Set-TextValue $ws.Range("E36") "  +0.18%  "
Set-TextValue $ws.Range("D37") "17.99"
Set-TextValue $ws.Range("E37") "  +1.20%  "
Set-TextValue $ws.Range("E38") "  -0.08%  "
Set-TextValue $ws.Range("D39") "3.84"
Set-TextValue $ws.Range("E39") "  +2.78%  "
Set-TextValue $ws.Range("D40") "0.838"
Set-TextValue $ws.Range("E40") "  +4.94%  "
Set-TextValue $ws.Range("D41") "36.86"
Set-TextValue $ws.Range("E41") "  +2.80%  "
Set-TextValue $ws.Range("E42") "  +0.34%  "
Set-TextValue $ws.Range("D43") "132.92"
Set-TextValue $ws.Range("E43") "  +2.65%  "
Set-TextValue $ws.Range("E44") "  +0.78%  "
Set-TextValue $ws.Range("E45") "  -0.82%  "
Set-TextValue $ws.Range("D46") "0.572"
Set-TextValue $ws.Range("E46") "  -0.79%  "
Set-TextValue $ws.Range("D47") "0.0914"
Set-TextValue $ws.Range("E47") "  +1.55%  "
Set-TextValue $ws.Range("D48") "251.18"
Set-TextValue $ws.Range("E48") "  -1.97%  "
Set-TextValue $ws.Range("E49") "  -0.04%  "
Set-TextValue $ws.Range("E50") "  +2.00%  "
Set-TextValue $ws.Range("D51") "17.09"
Set-TextValue $ws.Range("E51") "  +7.14%  "
